$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing rows (row 26 was last).
$newRows = @(
    @(42602.513611111113, "Bag", 53, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42602.516736111109, "Bag", 18, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42602.524085648147, "Bag", 63, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0)
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A: date/time serial value, formatted like the other date cells.
    $srcA = $ws.Cells.Item($r - 1, 1)
    $cellA = $ws.Cells.Item($r, 1)
    $srcA.Copy($cellA)
    $cellA.Value = $row[0]

    # Column B: shared string "Bag".
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Columns C..M: numeric values.
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
